# More complete reference data
# Fills in previously-blank Isolate/Year/Country/pubmedID cells on the
# "main refs" sheet (rows 3-12), and adds missing pubmedID values for
# rows 13-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main refs")

$xlPasteFormats = -4122

function Set-FormattedValue {
    param(
        [string]$TargetCell,
        [string]$StyleSourceCell,
        $Value
    )
    if ($StyleSourceCell) {
        $ws.Range($StyleSourceCell).Copy()
        $ws.Range($TargetCell).PasteSpecial($xlPasteFormats)
    }
    if ($null -ne $Value) {
        $ws.Range($TargetCell).Value = $Value
    }
}

# --- Row 3 (JN188292) ---
Set-FormattedValue "E3" "D3" "pZAC"
Set-FormattedValue "F3" "D3" 1981
Set-FormattedValue "G3" "D3" $null
Set-FormattedValue "J3" "D3" 23170185

# --- Row 4 (KU749402) ---
Set-FormattedValue "E4" "D4" "DEMA112UA014"
Set-FormattedValue "F4" "D4" 2012
Set-FormattedValue "G4" "D4" $null
Set-FormattedValue "J4" "D4" "NULL"

# --- Row 5 (KU168275) ---
Set-FormattedValue "E5" "D5" "LA21LeAn"
Set-FormattedValue "F5" "D5" 2003
Set-FormattedValue "G5" "D5" $null
Set-FormattedValue "J5" "D5" 26699702

# --- Row 6 (AF385935) ---
Set-FormattedValue "E6" "D6" "URTR35"
Set-FormattedValue "F6" "D6" 1999
Set-FormattedValue "G6" "D6" $null
Set-FormattedValue "J6" "D6" 15585101
$ws.Rows.Item(6).AutoFit()

# --- Row 7 (EF036536) ---
Set-FormattedValue "E7" "D7" "Fj061"
Set-FormattedValue "F7" "D7" 2006
Set-FormattedValue "J7" "D7" 17451347

# --- Row 8 (KC503852) ---
Set-FormattedValue "E8" "D8" "pXJDC6291-2-6"
Set-FormattedValue "F8" "D8" 2005
Set-FormattedValue "G8" "D8" "China"
Set-FormattedValue "J8" "D8" 24324545

# --- Row 9 (KU168273) ---
Set-FormattedValue "G9" "D9" "Democratic Republic of the Congo"
Set-FormattedValue "E9" "D9" "LA19KoSa"
Set-FormattedValue "F9" "D9" 2004
Set-FormattedValue "J9" "D9" 26699702

# --- Row 10 (AJ006022) ---
Set-FormattedValue "F10" "D10" 1995
Set-FormattedValue "G10" "D10" "Cameroon"
Set-FormattedValue "J10" "D10" 9734396

# --- Row 11 (KU168292) ---
Set-FormattedValue "E11" "D11" "LA49RBF189"
Set-FormattedValue "F11" "D11" 2005
Set-FormattedValue "G11" "D11" "Cameroon"
Set-FormattedValue "J11" "D11" 26699702

# --- Row 12 (GU111555) ---
Set-FormattedValue "E12" "D12" "RBF168"
Set-FormattedValue "F12" "D12" 2009
Set-FormattedValue "G12" "D12" "France"
Set-FormattedValue "J12" "D12" "NULL"

# Row 10's isolate is filled in last (matches original authoring order).
$ws.Range("E10").ClearFormats()
$ws.Range("E10").Value = "YBF30"

# --- Rows 13-17: add missing pubmedID values only ---
$ws.Range("J13").Value = 25733890
$ws.Range("J14").Value = 25733890
$ws.Range("J15").Value = 22505456
$ws.Range("J16").Value = 22505456
$ws.Range("J17").Value = 17494082
